$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update last-updated timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Cells.Item(2, 1).Value = "30 Oct 2025, 01:15 PM"

# --- Top Gainers sheet ---
$wsGainers = $wb.Worksheets.Item("Top Gainers")
$wsGainers.Cells.Item(2, 3).Value = 11.066
$wsGainers.Cells.Item(2, 4).Value = 16.383
$wsGainers.Cells.Item(2, 5).Value = 23.3371
$wsGainers.Cells.Item(3, 3).Value = 10.8174
$wsGainers.Cells.Item(3, 4).Value = 19.5954
$wsGainers.Cells.Item(3, 5).Value = 26.73
$wsGainers.Cells.Item(4, 3).Value = 10.4629
$wsGainers.Cells.Item(4, 4).Value = 10.6693
$wsGainers.Cells.Item(4, 5).Value = 24.73
$wsGainers.Cells.Item(5, 3).Value = 10.3896
$wsGainers.Cells.Item(5, 4).Value = 7.438
$wsGainers.Cells.Item(5, 5).Value = -7.4152
$wsGainers.Cells.Item(6, 2).Value = "SAGILITY"
$wsGainers.Cells.Item(6, 3).Value = 9.8093
$wsGainers.Cells.Item(6, 4).Value = 17.5505
$wsGainers.Cells.Item(6, 5).Value = 30.9423
$wsGainers.Cells.Item(7, 2).Value = "UNIPARTS"
$wsGainers.Cells.Item(7, 3).Value = 9.552300000000001
$wsGainers.Cells.Item(7, 4).Value = 11.8682
$wsGainers.Cells.Item(7, 5).Value = 28.3393
$wsGainers.Cells.Item(8, 3).Value = 7.5707
$wsGainers.Cells.Item(8, 4).Value = 12.7727
$wsGainers.Cells.Item(8, 5).Value = 14.938
$wsGainers.Cells.Item(9, 2).Value = "MCLOUD"
$wsGainers.Cells.Item(9, 3).Value = 6.4337
$wsGainers.Cells.Item(9, 4).Value = 5.4136
$wsGainers.Cells.Item(9, 5).Value = -22.8801
$wsGainers.Cells.Item(10, 3).Value = 6.2082
$wsGainers.Cells.Item(10, 4).Value = 10.4865
$wsGainers.Cells.Item(10, 5).Value = 11.5823
$wsGainers.Cells.Item(11, 2).Value = "MARINE"
$wsGainers.Cells.Item(11, 3).Value = 6.0897
$wsGainers.Cells.Item(11, 4).Value = 2.627
$wsGainers.Cells.Item(11, 5).Value = 15.1593
$wsGainers.Cells.Item(12, 2).Value = "VENKEYS"
$wsGainers.Cells.Item(12, 3).Value = 5.7234
$wsGainers.Cells.Item(12, 4).Value = 6.3547
$wsGainers.Cells.Item(12, 5).Value = 4.0414
$wsGainers.Cells.Item(13, 3).Value = 5.7235
$wsGainers.Cells.Item(13, 4).Value = 11.6139
$wsGainers.Cells.Item(13, 5).Value = 13.7125
$wsGainers.Cells.Item(14, 3).Value = 5.661
$wsGainers.Cells.Item(14, 4).Value = 4.6248
$wsGainers.Cells.Item(14, 5).Value = -1.583
$wsGainers.Cells.Item(15, 2).Value = "PDSL"
$wsGainers.Cells.Item(15, 3).Value = 5.5796
$wsGainers.Cells.Item(15, 4).Value = 8.9246
$wsGainers.Cells.Item(15, 5).Value = 15.0872
$wsGainers.Cells.Item(16, 2).Value = "INDIACEM"
$wsGainers.Cells.Item(16, 3).Value = 5.4044
$wsGainers.Cells.Item(16, 4).Value = 5.6879
$wsGainers.Cells.Item(16, 5).Value = 7.4219
$wsGainers.Cells.Item(17, 2).Value = "POLICYBZR"
$wsGainers.Cells.Item(17, 3).Value = 5.3905
$wsGainers.Cells.Item(17, 4).Value = 7.7475
$wsGainers.Cells.Item(17, 5).Value = 6.7156
$wsGainers.Cells.Item(21, 2).Value = "HIRECT"
$wsGainers.Cells.Item(21, 3).Value = 4.9834
$wsGainers.Cells.Item(21, 4).Value = 12.4953
$wsGainers.Cells.Item(21, 5).Value = 10.7678
$wsGainers.Cells.Item(22, 3).Value = 4.937
$wsGainers.Cells.Item(22, 4).Value = 3.5501
$wsGainers.Cells.Item(22, 5).Value = 8.8531
$wsGainers.Cells.Item(23, 2).Value = "VIMTALABS"
$wsGainers.Cells.Item(23, 3).Value = 4.9246
$wsGainers.Cells.Item(23, 4).Value = 5.1042
$wsGainers.Cells.Item(23, 5).Value = -0.0283
$wsGainers.Cells.Item(24, 2).Value = "EUROPRATIK"
$wsGainers.Cells.Item(24, 3).Value = 4.825
$wsGainers.Cells.Item(24, 4).Value = 10.7264
$wsGainers.Cells.Item(24, 5).Value = 27.7646
$wsGainers.Cells.Item(25, 2).Value = "RAMASTEEL"
$wsGainers.Cells.Item(25, 3).Value = 4.6185
$wsGainers.Cells.Item(25, 4).Value = 4.5135
$wsGainers.Cells.Item(25, 5).Value = 6.11
$wsGainers.Cells.Item(26, 3).Value = 4.4949
$wsGainers.Cells.Item(26, 4).Value = 11.0332
$wsGainers.Cells.Item(26, 5).Value = 7.5362
$wsGainers.Cells.Item(27, 2).Value = "BAJAJHCARE"
$wsGainers.Cells.Item(27, 3).Value = 4.237
$wsGainers.Cells.Item(27, 4).Value = 4.7848
$wsGainers.Cells.Item(27, 5).Value = -1.5027
$wsGainers.Cells.Item(28, 2).Value = "ALICON"
$wsGainers.Cells.Item(28, 3).Value = 4.2351
$wsGainers.Cells.Item(28, 4).Value = 10.4143
$wsGainers.Cells.Item(28, 5).Value = 15.8988
$wsGainers.Cells.Item(29, 3).Value = 4.1861
$wsGainers.Cells.Item(29, 4).Value = 8.3344
$wsGainers.Cells.Item(29, 5).Value = 32.4133
$wsGainers.Cells.Item(30, 3).Value = 4.1645
$wsGainers.Cells.Item(30, 4).Value = 5.2659
$wsGainers.Cells.Item(30, 5).Value = 6.3571
$wsGainers.Cells.Item(31, 2).Value = "INDORAMA"
$wsGainers.Cells.Item(31, 3).Value = 3.9359
$wsGainers.Cells.Item(31, 4).Value = 6.4338
$wsGainers.Cells.Item(31, 5).Value = 17.9682
$wsGainers.Cells.Item(33, 3).Value = 3.8673
$wsGainers.Cells.Item(33, 4).Value = 12.2266
$wsGainers.Cells.Item(33, 5).Value = 11.7616
$wsGainers.Cells.Item(34, 2).Value = "DEEDEV"
$wsGainers.Cells.Item(34, 3).Value = 3.8672
$wsGainers.Cells.Item(34, 4).Value = -3.0429
$wsGainers.Cells.Item(34, 5).Value = -3.8425
$wsGainers.Cells.Item(35, 2).Value = "SKYGOLD"
$wsGainers.Cells.Item(35, 3).Value = 3.8478
$wsGainers.Cells.Item(35, 4).Value = -0.7112000000000001
$wsGainers.Cells.Item(35, 5).Value = 37.8917
$wsGainers.Cells.Item(36, 2).Value = "PFOCUS"
$wsGainers.Cells.Item(36, 3).Value = 3.7112
$wsGainers.Cells.Item(36, 4).Value = 0.986
$wsGainers.Cells.Item(36, 5).Value = 2.4497
$wsGainers.Cells.Item(37, 2).Value = "SHANTIGOLD"
$wsGainers.Cells.Item(37, 3).Value = 3.6832
$wsGainers.Cells.Item(37, 4).Value = 10.9809
$wsGainers.Cells.Item(37, 5).Value = 3.5766
$wsGainers.Cells.Item(38, 3).Value = 3.634
$wsGainers.Cells.Item(38, 4).Value = 6.9012
$wsGainers.Cells.Item(38, 5).Value = -0.3405
$wsGainers.Cells.Item(39, 2).Value = "SAPPHIRE"
$wsGainers.Cells.Item(39, 3).Value = 3.6198
$wsGainers.Cells.Item(39, 4).Value = 5.4469
$wsGainers.Cells.Item(39, 5).Value = 2.791
$wsGainers.Cells.Item(40, 2).Value = "CANBK"
$wsGainers.Cells.Item(40, 3).Value = 3.58
$wsGainers.Cells.Item(40, 4).Value = 6.1098
$wsGainers.Cells.Item(40, 5).Value = 7.8167
$wsGainers.Cells.Item(41, 2).Value = "GRMOVER"
$wsGainers.Cells.Item(41, 3).Value = 3.5257
$wsGainers.Cells.Item(41, 4).Value = 3.701
$wsGainers.Cells.Item(41, 5).Value = 19.465
$wsGainers.Cells.Item(42, 2).Value = "GMMPFAUDLR"
$wsGainers.Cells.Item(42, 3).Value = 3.521
$wsGainers.Cells.Item(42, 4).Value = 8.007400000000001
$wsGainers.Cells.Item(42, 5).Value = 20.4114
$wsGainers.Cells.Item(43, 2).Value = "AHLUCONT"
$wsGainers.Cells.Item(43, 3).Value = 3.5124
$wsGainers.Cells.Item(43, 4).Value = 2.3145
$wsGainers.Cells.Item(43, 5).Value = -4.9437
$wsGainers.Cells.Item(44, 2).Value = "REDTAPE"
$wsGainers.Cells.Item(44, 3).Value = 3.4967
$wsGainers.Cells.Item(44, 4).Value = 3.3963
$wsGainers.Cells.Item(44, 5).Value = -3.5242
$wsGainers.Cells.Item(45, 2).Value = "SUNDROP"
$wsGainers.Cells.Item(45, 3).Value = 3.3957
$wsGainers.Cells.Item(45, 4).Value = 3.288
$wsGainers.Cells.Item(45, 5).Value = 1.3804
$wsGainers.Cells.Item(46, 2).Value = "FIVESTAR"
$wsGainers.Cells.Item(46, 3).Value = 3.3386
$wsGainers.Cells.Item(46, 4).Value = 16.4339
$wsGainers.Cells.Item(46, 5).Value = 16.5207
$wsGainers.Cells.Item(47, 2).Value = "JKTYRE"
$wsGainers.Cells.Item(47, 3).Value = 3.3361
$wsGainers.Cells.Item(47, 4).Value = 6.3964
$wsGainers.Cells.Item(47, 5).Value = 22.55
$wsGainers.Cells.Item(48, 2).Value = "VSTIND"
$wsGainers.Cells.Item(48, 3).Value = 3.3262
$wsGainers.Cells.Item(48, 4).Value = 3.8282
$wsGainers.Cells.Item(48, 5).Value = 3.3662
$wsGainers.Cells.Item(49, 2).Value = "BLISSGVS"
$wsGainers.Cells.Item(49, 3).Value = 3.2752
$wsGainers.Cells.Item(49, 4).Value = 2.6273
$wsGainers.Cells.Item(49, 5).Value = 2.9975
$wsGainers.Cells.Item(50, 2).Value = "BLUEDART"
$wsGainers.Cells.Item(50, 3).Value = 3.1954
$wsGainers.Cells.Item(50, 4).Value = 21.8579
$wsGainers.Cells.Item(50, 5).Value = 18.9199
$wsGainers.Cells.Item(51, 2).Value = "NEULANDLAB"
$wsGainers.Cells.Item(51, 3).Value = 3.1871
$wsGainers.Cells.Item(51, 4).Value = -1.1872
$wsGainers.Cells.Item(51, 5).Value = 8.846299999999999
$wsGainers.Cells.Item(52, 3).Value = 3.1578
$wsGainers.Cells.Item(52, 4).Value = 9.941599999999999
$wsGainers.Cells.Item(52, 5).Value = -1.5758
$wsGainers.Cells.Item(54, 3).Value = 3.0303
$wsGainers.Cells.Item(54, 4).Value = 1.6442
$wsGainers.Cells.Item(54, 5).Value = 0.9801
$wsGainers.Cells.Item(55, 3).Value = 2.9356
$wsGainers.Cells.Item(55, 4).Value = 12.931
$wsGainers.Cells.Item(55, 5).Value = 23.5786
$wsGainers.Cells.Item(56, 2).Value = "ASHOKA"
$wsGainers.Cells.Item(56, 3).Value = 2.8565
$wsGainers.Cells.Item(56, 4).Value = 4.3772
$wsGainers.Cells.Item(56, 5).Value = 7.0044
$wsGainers.Cells.Item(57, 2).Value = "VOLTAMP"
$wsGainers.Cells.Item(57, 3).Value = 2.83
$wsGainers.Cells.Item(57, 4).Value = 2.7089
$wsGainers.Cells.Item(57, 5).Value = 2.4463
$wsGainers.Cells.Item(58, 3).Value = 2.82
$wsGainers.Cells.Item(58, 4).Value = 7.2046
$wsGainers.Cells.Item(58, 5).Value = -0.1444
$wsGainers.Cells.Item(59, 3).Value = 2.7938
$wsGainers.Cells.Item(59, 4).Value = 16.8568
$wsGainers.Cells.Item(59, 5).Value = 23.2815
$wsGainers.Cells.Item(60, 2).Value = "WESTLIFE"
$wsGainers.Cells.Item(60, 3).Value = 2.7496
$wsGainers.Cells.Item(60, 4).Value = 2.7059
$wsGainers.Cells.Item(60, 5).Value = -12.2309
$wsGainers.Cells.Item(61, 2).Value = "BPCL"
$wsGainers.Cells.Item(61, 3).Value = 2.7291
$wsGainers.Cells.Item(61, 4).Value = 8.216100000000001
$wsGainers.Cells.Item(61, 5).Value = 5.2849
$wsGainers.Cells.Item(62, 2).Value = "BGRENERGY"
$wsGainers.Cells.Item(62, 3).Value = 2.7202
$wsGainers.Cells.Item(62, 4).Value = -6.5421
$wsGainers.Cells.Item(62, 5).Value = 73.9896
$wsGainers.Cells.Item(63, 2).Value = "INOXGREEN"
$wsGainers.Cells.Item(63, 3).Value = 2.6165
$wsGainers.Cells.Item(63, 4).Value = 10.6171
$wsGainers.Cells.Item(63, 5).Value = 33.9715
$wsGainers.Cells.Item(64, 2).Value = "GANESHCP"
$wsGainers.Cells.Item(64, 3).Value = 2.585
$wsGainers.Cells.Item(64, 4).Value = 2.0609
$wsGainers.Cells.Item(64, 5).Value = 1.6135
$wsGainers.Cells.Item(65, 2).Value = "BLS"
$wsGainers.Cells.Item(65, 3).Value = 2.5255
$wsGainers.Cells.Item(65, 4).Value = -0.5028
$wsGainers.Cells.Item(65, 5).Value = -1.753
$wsGainers.Cells.Item(66, 2).Value = "MFSL"
$wsGainers.Cells.Item(66, 3).Value = 2.5151
$wsGainers.Cells.Item(66, 4).Value = 2.5692
$wsGainers.Cells.Item(66, 5).Value = -1.1993
$wsGainers.Cells.Item(67, 2).Value = "JKLAKSHMI"
$wsGainers.Cells.Item(67, 3).Value = 2.4909
$wsGainers.Cells.Item(67, 4).Value = 4.4781
$wsGainers.Cells.Item(67, 5).Value = 1.4987
$wsGainers.Cells.Item(68, 2).Value = "REFEX"
$wsGainers.Cells.Item(68, 3).Value = 2.4542
$wsGainers.Cells.Item(68, 4).Value = -0.11
$wsGainers.Cells.Item(68, 5).Value = 1.8651
$wsGainers.Cells.Item(69, 2).Value = "OIL"
$wsGainers.Cells.Item(69, 3).Value = 2.4495
$wsGainers.Cells.Item(69, 4).Value = 2.6937
$wsGainers.Cells.Item(69, 5).Value = 4.1083
$wsGainers.Cells.Item(70, 2).Value = "IIFL"
$wsGainers.Cells.Item(70, 3).Value = 2.4474
$wsGainers.Cells.Item(70, 4).Value = 9.3916
$wsGainers.Cells.Item(70, 5).Value = 18.5661
$wsGainers.Cells.Item(71, 2).Value = "CIFL"
$wsGainers.Cells.Item(71, 3).Value = 2.4426
$wsGainers.Cells.Item(71, 4).Value = 1.9977
$wsGainers.Cells.Item(71, 5).Value = 1.9387
$wsGainers.Cells.Item(72, 2).Value = "FEDFINA"
$wsGainers.Cells.Item(72, 3).Value = 2.4139
$wsGainers.Cells.Item(72, 4).Value = 3.5149
$wsGainers.Cells.Item(72, 5).Value = -5.2131
$wsGainers.Cells.Item(73, 2).Value = "OBEROIRLTY"
$wsGainers.Cells.Item(73, 3).Value = 2.4104
$wsGainers.Cells.Item(73, 4).Value = 3.2237
$wsGainers.Cells.Item(73, 5).Value = 10.9095
$wsGainers.Cells.Item(75, 2).Value = "CENTUM"
$wsGainers.Cells.Item(75, 3).Value = 2.3583
$wsGainers.Cells.Item(75, 4).Value = 3.2128
$wsGainers.Cells.Item(75, 5).Value = -2.2594
$wsGainers.Cells.Item(76, 2).Value = "SDBL"
$wsGainers.Cells.Item(76, 3).Value = 2.3517
$wsGainers.Cells.Item(76, 4).Value = 0.8999
$wsGainers.Cells.Item(76, 5).Value = 6.4783

# --- Top Losers sheet ---
$wsLosers = $wb.Worksheets.Item("Top Losers")
$wsLosers.Cells.Item(2, 3).Value = -17.2902
$wsLosers.Cells.Item(2, 4).Value = -16.01
$wsLosers.Cells.Item(2, 5).Value = 1.282
$wsLosers.Cells.Item(3, 3).Value = -9.764699999999999
$wsLosers.Cells.Item(3, 4).Value = -6.3359
$wsLosers.Cells.Item(3, 5).Value = 5.8707
$wsLosers.Cells.Item(4, 3).Value = -6.8376
$wsLosers.Cells.Item(4, 4).Value = -9.355499999999999
$wsLosers.Cells.Item(4, 5).Value = 7.2571
$wsLosers.Cells.Item(10, 3).Value = -4.8938
$wsLosers.Cells.Item(10, 4).Value = -3.2341
$wsLosers.Cells.Item(10, 5).Value = 19.8335
$wsLosers.Cells.Item(11, 2).Value = "RAMCOSYS"
$wsLosers.Cells.Item(11, 3).Value = -4.5458
$wsLosers.Cells.Item(11, 4).Value = 5.1452
$wsLosers.Cells.Item(11, 5).Value = 23.5068
$wsLosers.Cells.Item(12, 2).Value = "NSLNISP"
$wsLosers.Cells.Item(12, 3).Value = -4.3932
$wsLosers.Cells.Item(12, 4).Value = -3.1468
$wsLosers.Cells.Item(12, 5).Value = -3.9456
$wsLosers.Cells.Item(13, 2).Value = "RAJRATAN"
$wsLosers.Cells.Item(13, 3).Value = -4.3625
$wsLosers.Cells.Item(13, 4).Value = -3.5791
$wsLosers.Cells.Item(13, 5).Value = 21.2712
$wsLosers.Cells.Item(14, 2).Value = "SOUTHBANK"
$wsLosers.Cells.Item(14, 3).Value = -4.3212
$wsLosers.Cells.Item(14, 4).Value = -1.8106
$wsLosers.Cells.Item(14, 5).Value = 29.2573
$wsLosers.Cells.Item(15, 2).Value = "LICHSGFIN"
$wsLosers.Cells.Item(15, 3).Value = -4.2969
$wsLosers.Cells.Item(15, 4).Value = -1.9931
$wsLosers.Cells.Item(15, 5).Value = 0.5221
$wsLosers.Cells.Item(16, 2).Value = "SARDAEN"
$wsLosers.Cells.Item(16, 3).Value = -3.9277
$wsLosers.Cells.Item(16, 4).Value = -0.4266
$wsLosers.Cells.Item(16, 5).Value = -0.4635
$wsLosers.Cells.Item(17, 2).Value = "IDEAFORGE"
$wsLosers.Cells.Item(17, 3).Value = -3.9051
$wsLosers.Cells.Item(17, 4).Value = -3.1095
$wsLosers.Cells.Item(17, 5).Value = -4.8389
$wsLosers.Cells.Item(18, 2).Value = "YATRA"
$wsLosers.Cells.Item(18, 3).Value = -3.8095
$wsLosers.Cells.Item(18, 4).Value = -6.5465
$wsLosers.Cells.Item(18, 5).Value = 3.2808
$wsLosers.Cells.Item(19, 2).Value = "KHAICHEM"
$wsLosers.Cells.Item(19, 3).Value = -3.6538
$wsLosers.Cells.Item(19, 4).Value = -10.1288
$wsLosers.Cells.Item(19, 5).Value = -8.1295
$wsLosers.Cells.Item(20, 3).Value = -3.5974
$wsLosers.Cells.Item(20, 4).Value = -6.0533
$wsLosers.Cells.Item(20, 5).Value = -1.4546
$wsLosers.Cells.Item(21, 2).Value = "IEX"
$wsLosers.Cells.Item(21, 3).Value = -3.5786
$wsLosers.Cells.Item(21, 4).Value = -2.523
$wsLosers.Cells.Item(21, 5).Value = 2.9815
$wsLosers.Cells.Item(22, 2).Value = "MEGASOFT"
$wsLosers.Cells.Item(22, 3).Value = -3.5611
$wsLosers.Cells.Item(22, 4).Value = 11.6365
$wsLosers.Cells.Item(22, 5).Value = 28.772
$wsLosers.Cells.Item(23, 2).Value = "INDUSTOWER"
$wsLosers.Cells.Item(23, 3).Value = -3.5166
$wsLosers.Cells.Item(23, 4).Value = 1.6872
$wsLosers.Cells.Item(23, 5).Value = 7.2178
$wsLosers.Cells.Item(24, 2).Value = "TVSHLTD"
$wsLosers.Cells.Item(24, 3).Value = -3.4813
$wsLosers.Cells.Item(24, 4).Value = -2.2385
$wsLosers.Cells.Item(24, 5).Value = 16.0266
$wsLosers.Cells.Item(25, 2).Value = "VBL"
$wsLosers.Cells.Item(25, 3).Value = -3.4009
$wsLosers.Cells.Item(25, 4).Value = 3.7278
$wsLosers.Cells.Item(25, 5).Value = 7.8657
$wsLosers.Cells.Item(26, 2).Value = "APOLLOPIPE"
$wsLosers.Cells.Item(26, 3).Value = -3.3882
$wsLosers.Cells.Item(26, 4).Value = -5.0729
$wsLosers.Cells.Item(26, 5).Value = -10.0901
$wsLosers.Cells.Item(27, 2).Value = "SANDHAR"
$wsLosers.Cells.Item(27, 3).Value = -3.3167
$wsLosers.Cells.Item(27, 4).Value = 0.3741
$wsLosers.Cells.Item(27, 5).Value = 17.6239
$wsLosers.Cells.Item(28, 2).Value = "VGUARD"
$wsLosers.Cells.Item(28, 3).Value = -3.2545
$wsLosers.Cells.Item(28, 4).Value = -0.5642
$wsLosers.Cells.Item(28, 5).Value = -1.4381
$wsLosers.Cells.Item(29, 2).Value = "UBL"
$wsLosers.Cells.Item(29, 3).Value = -3.2277
$wsLosers.Cells.Item(29, 4).Value = -2.6352
$wsLosers.Cells.Item(29, 5).Value = -1.1051
$wsLosers.Cells.Item(30, 2).Value = "SAIL"
$wsLosers.Cells.Item(30, 3).Value = -3.095
$wsLosers.Cells.Item(30, 4).Value = 5.2062
$wsLosers.Cells.Item(30, 5).Value = 1.279
$wsLosers.Cells.Item(31, 2).Value = "TVSELECT"
$wsLosers.Cells.Item(31, 3).Value = -3.0894
$wsLosers.Cells.Item(31, 4).Value = -4.0331
$wsLosers.Cells.Item(31, 5).Value = -5.9936
$wsLosers.Cells.Item(32, 2).Value = "FILATEX"
$wsLosers.Cells.Item(32, 3).Value = -3.059
$wsLosers.Cells.Item(32, 4).Value = 6.9007
$wsLosers.Cells.Item(32, 5).Value = 22.1483
$wsLosers.Cells.Item(33, 2).Value = "GOKULAGRO"
$wsLosers.Cells.Item(33, 3).Value = -3.0312
$wsLosers.Cells.Item(33, 4).Value = 4.5314
$wsLosers.Cells.Item(33, 5).Value = -13.8597
$wsLosers.Cells.Item(34, 2).Value = "IDBI"
$wsLosers.Cells.Item(34, 3).Value = -2.983
$wsLosers.Cells.Item(34, 4).Value = 5.2592
$wsLosers.Cells.Item(34, 5).Value = 8.1492
$wsLosers.Cells.Item(35, 2).Value = "DREDGECORP"
$wsLosers.Cells.Item(35, 3).Value = -2.9569
$wsLosers.Cells.Item(35, 4).Value = 18.1621
$wsLosers.Cells.Item(35, 5).Value = 18.9346
$wsLosers.Cells.Item(36, 2).Value = "JSL"
$wsLosers.Cells.Item(36, 3).Value = -2.9489
$wsLosers.Cells.Item(36, 4).Value = -2.7492
$wsLosers.Cells.Item(36, 5).Value = 5.7269
$wsLosers.Cells.Item(37, 2).Value = "LXCHEM"
$wsLosers.Cells.Item(37, 3).Value = -2.9431
$wsLosers.Cells.Item(37, 4).Value = -3.3805
$wsLosers.Cells.Item(37, 5).Value = -4.5496
$wsLosers.Cells.Item(41, 2).Value = "ARIHANTCAP"
$wsLosers.Cells.Item(41, 3).Value = -2.796
$wsLosers.Cells.Item(41, 4).Value = 4.1311
$wsLosers.Cells.Item(41, 5).Value = -4.648
$wsLosers.Cells.Item(42, 2).Value = "CGCL"
$wsLosers.Cells.Item(42, 3).Value = -2.7884
$wsLosers.Cells.Item(42, 4).Value = -0.8343
$wsLosers.Cells.Item(42, 5).Value = 9.758100000000001
$wsLosers.Cells.Item(43, 2).Value = "FABTECH"
$wsLosers.Cells.Item(43, 3).Value = -2.745
$wsLosers.Cells.Item(43, 4).Value = 12.4494
$wsLosers.Cells.Item(43, 5).Value = "N/A"
$wsLosers.Cells.Item(44, 3).Value = -2.7256
$wsLosers.Cells.Item(44, 4).Value = -7.3342
$wsLosers.Cells.Item(44, 5).Value = -4.7619
$wsLosers.Cells.Item(46, 2).Value = "BCG"
$wsLosers.Cells.Item(46, 3).Value = -2.7027
$wsLosers.Cells.Item(46, 4).Value = 2.3186
$wsLosers.Cells.Item(46, 5).Value = -1.5119
$wsLosers.Cells.Item(47, 2).Value = "BHARATWIRE"
$wsLosers.Cells.Item(47, 3).Value = -2.6745
$wsLosers.Cells.Item(47, 4).Value = 19.5484
$wsLosers.Cells.Item(47, 5).Value = 20.5843
$wsLosers.Cells.Item(48, 2).Value = "HFCL"
$wsLosers.Cells.Item(48, 3).Value = -2.6572
$wsLosers.Cells.Item(48, 4).Value = -3.3959
$wsLosers.Cells.Item(48, 5).Value = 3.0461
$wsLosers.Cells.Item(49, 2).Value = "GRWRHITECH"
$wsLosers.Cells.Item(49, 3).Value = -2.6221
$wsLosers.Cells.Item(49, 4).Value = -6.2062
$wsLosers.Cells.Item(49, 5).Value = 18.4864
$wsLosers.Cells.Item(50, 2).Value = "CAMLINFINE"
$wsLosers.Cells.Item(50, 3).Value = -2.6177
$wsLosers.Cells.Item(50, 4).Value = 0.1635
$wsLosers.Cells.Item(50, 5).Value = 0.4252
$wsLosers.Cells.Item(51, 2).Value = "STYL"
$wsLosers.Cells.Item(51, 3).Value = -2.6172
$wsLosers.Cells.Item(51, 4).Value = -5.7995
$wsLosers.Cells.Item(51, 5).Value = -11.3744
$wsLosers.Cells.Item(53, 2).Value = "BHARTIHEXA"
$wsLosers.Cells.Item(53, 3).Value = -2.5754
$wsLosers.Cells.Item(53, 4).Value = 4.3298
$wsLosers.Cells.Item(53, 5).Value = 12.3629
$wsLosers.Cells.Item(54, 3).Value = -2.5557
$wsLosers.Cells.Item(54, 4).Value = -8.4131
$wsLosers.Cells.Item(54, 5).Value = 21.6196
$wsLosers.Cells.Item(56, 2).Value = "POCL"
$wsLosers.Cells.Item(56, 3).Value = -2.4488
$wsLosers.Cells.Item(56, 4).Value = 2.7169
$wsLosers.Cells.Item(56, 5).Value = 23.2477
$wsLosers.Cells.Item(57, 2).Value = "KFINTECH"
$wsLosers.Cells.Item(57, 3).Value = -2.3909
$wsLosers.Cells.Item(57, 4).Value = -4.0978
$wsLosers.Cells.Item(57, 5).Value = 4.8113
$wsLosers.Cells.Item(58, 3).Value = -2.3702
$wsLosers.Cells.Item(58, 4).Value = -2.3048
$wsLosers.Cells.Item(59, 2).Value = "NUVAMA"
$wsLosers.Cells.Item(59, 3).Value = -2.3445
$wsLosers.Cells.Item(59, 4).Value = 0.1254
$wsLosers.Cells.Item(59, 5).Value = 13.9572
$wsLosers.Cells.Item(60, 2).Value = "AEGISLOG"
$wsLosers.Cells.Item(60, 3).Value = -2.3238
$wsLosers.Cells.Item(60, 4).Value = -1.9806
$wsLosers.Cells.Item(60, 5).Value = 0.1763
$wsLosers.Cells.Item(61, 2).Value = "SOLARWORLD"
$wsLosers.Cells.Item(61, 3).Value = -2.3197
$wsLosers.Cells.Item(61, 4).Value = 6.4916
$wsLosers.Cells.Item(61, 5).Value = 2.1793
$wsLosers.Cells.Item(62, 2).Value = "PROSTARM"
$wsLosers.Cells.Item(62, 3).Value = -2.3139
$wsLosers.Cells.Item(62, 4).Value = -1.6574
$wsLosers.Cells.Item(62, 5).Value = -10.3525
$wsLosers.Cells.Item(63, 2).Value = "BEPL"
$wsLosers.Cells.Item(63, 3).Value = -2.3096
$wsLosers.Cells.Item(63, 4).Value = 0.5361
$wsLosers.Cells.Item(63, 5).Value = -2.0786
$wsLosers.Cells.Item(64, 2).Value = "QUESS"
$wsLosers.Cells.Item(64, 3).Value = -2.308
$wsLosers.Cells.Item(64, 4).Value = 4.8697
$wsLosers.Cells.Item(64, 5).Value = -3.5742
$wsLosers.Cells.Item(65, 2).Value = "VIPIND"
$wsLosers.Cells.Item(65, 3).Value = -2.3015
$wsLosers.Cells.Item(65, 4).Value = -4.0295
$wsLosers.Cells.Item(65, 5).Value = -1.5146
$wsLosers.Cells.Item(66, 2).Value = "SURYAROSNI"
$wsLosers.Cells.Item(66, 3).Value = -2.2716
$wsLosers.Cells.Item(66, 4).Value = 8.855700000000001
$wsLosers.Cells.Item(66, 5).Value = 0.6811
$wsLosers.Cells.Item(67, 2).Value = "ARMANFIN"
$wsLosers.Cells.Item(67, 3).Value = -2.2632
$wsLosers.Cells.Item(67, 4).Value = -3.2436
$wsLosers.Cells.Item(67, 5).Value = 11.5435
$wsLosers.Cells.Item(68, 2).Value = "NEWGEN"
$wsLosers.Cells.Item(68, 3).Value = -2.2596
$wsLosers.Cells.Item(68, 4).Value = 9.017099999999999
$wsLosers.Cells.Item(68, 5).Value = 9.3908
$wsLosers.Cells.Item(69, 2).Value = "63MOONS"
$wsLosers.Cells.Item(69, 3).Value = -2.2501
$wsLosers.Cells.Item(69, 4).Value = 1.145
$wsLosers.Cells.Item(69, 5).Value = -5.0803
$wsLosers.Cells.Item(70, 2).Value = "MSPL"
$wsLosers.Cells.Item(70, 3).Value = -2.2449
$wsLosers.Cells.Item(70, 4).Value = -1.0915
$wsLosers.Cells.Item(70, 5).Value = -8.137
$wsLosers.Cells.Item(71, 2).Value = "CHAMBLFERT"
$wsLosers.Cells.Item(71, 3).Value = -2.2406
$wsLosers.Cells.Item(71, 4).Value = -0.7582
$wsLosers.Cells.Item(71, 5).Value = -5.1787
$wsLosers.Cells.Item(72, 2).Value = "GABRIEL"
$wsLosers.Cells.Item(72, 3).Value = -2.226
$wsLosers.Cells.Item(72, 4).Value = 1.9078
$wsLosers.Cells.Item(72, 5).Value = 6.7039
$wsLosers.Cells.Item(73, 2).Value = "COSMOFIRST"
$wsLosers.Cells.Item(73, 3).Value = -2.2188
$wsLosers.Cells.Item(73, 4).Value = -1.1593
$wsLosers.Cells.Item(73, 5).Value = 0.0354
$wsLosers.Cells.Item(74, 2).Value = "TMB"
$wsLosers.Cells.Item(74, 3).Value = -2.2155
$wsLosers.Cells.Item(74, 4).Value = 7.3642
$wsLosers.Cells.Item(74, 5).Value = 14.5396
$wsLosers.Cells.Item(75, 2).Value = "SINDHUTRAD"
$wsLosers.Cells.Item(75, 3).Value = -2.2014
$wsLosers.Cells.Item(75, 4).Value = -1.1129
$wsLosers.Cells.Item(75, 5).Value = -15.0796

# --- 1 Month Performance sheet ---
$wsPerf = $wb.Worksheets.Item("1 Month Performance")
$wsPerf.Cells.Item(4, 3).Value = 78.2526
$wsPerf.Cells.Item(5, 3).Value = 64.66160000000001
$wsPerf.Cells.Item(6, 3).Value = 61.353
$wsPerf.Cells.Item(8, 3).Value = 53.7404
$wsPerf.Cells.Item(9, 3).Value = 50.8068
$wsPerf.Cells.Item(10, 3).Value = 45.9663
$wsPerf.Cells.Item(11, 3).Value = 41.8364
$wsPerf.Cells.Item(14, 3).Value = 38.9267
$wsPerf.Cells.Item(15, 3).Value = 38.8627
$wsPerf.Cells.Item(16, 3).Value = 37.371
$wsPerf.Cells.Item(18, 3).Value = 34.6311
$wsPerf.Cells.Item(19, 3).Value = 34.2662
$wsPerf.Cells.Item(20, 3).Value = 33.268
$wsPerf.Cells.Item(21, 3).Value = 33.2011
$wsPerf.Cells.Item(22, 2).Value = "INDORAMA"
$wsPerf.Cells.Item(22, 3).Value = 32.9976
$wsPerf.Cells.Item(23, 3).Value = 32.9734
$wsPerf.Cells.Item(24, 2).Value = "MEGASOFT"
$wsPerf.Cells.Item(24, 3).Value = 32.6002
$wsPerf.Cells.Item(25, 3).Value = 30.577
$wsPerf.Cells.Item(26, 3).Value = 29.5706
$wsPerf.Cells.Item(27, 2).Value = "ONMOBILE"
$wsPerf.Cells.Item(27, 3).Value = 29.0018
$wsPerf.Cells.Item(28, 2).Value = "TARACHAND"
$wsPerf.Cells.Item(28, 3).Value = 28.9527
$wsPerf.Cells.Item(30, 2).Value = "MRPL"
$wsPerf.Cells.Item(30, 3).Value = 28.3335
$wsPerf.Cells.Item(31, 2).Value = "ADANIPOWER"
$wsPerf.Cells.Item(31, 3).Value = 28.3047
$wsPerf.Cells.Item(33, 3).Value = 25.7691
$wsPerf.Cells.Item(34, 3).Value = 25.6733
$wsPerf.Cells.Item(35, 3).Value = 25.5424
$wsPerf.Cells.Item(38, 3).Value = 25.2804
$wsPerf.Cells.Item(39, 3).Value = 24.8826
$wsPerf.Cells.Item(40, 2).Value = "MINDTECK"
$wsPerf.Cells.Item(40, 3).Value = 24.5083
$wsPerf.Cells.Item(41, 3).Value = 24.4742
$wsPerf.Cells.Item(42, 2).Value = "MARINE"
$wsPerf.Cells.Item(42, 3).Value = 24.2365
$wsPerf.Cells.Item(43, 3).Value = 24.1959
$wsPerf.Cells.Item(45, 3).Value = 22.8539
$wsPerf.Cells.Item(46, 2).Value = "DCBBANK"
$wsPerf.Cells.Item(46, 3).Value = 22.7476
$wsPerf.Cells.Item(49, 2).Value = "RAMCOSYS"
$wsPerf.Cells.Item(49, 3).Value = 22.5023
$wsPerf.Cells.Item(50, 3).Value = 22.4857
$wsPerf.Cells.Item(52, 3).Value = 21.6798
$wsPerf.Cells.Item(53, 3).Value = 21.665
$wsPerf.Cells.Item(54, 2).Value = "GUJTHEM"
$wsPerf.Cells.Item(54, 3).Value = 21.4951
$wsPerf.Cells.Item(55, 2).Value = "SCI"
$wsPerf.Cells.Item(55, 3).Value = 21.3462
$wsPerf.Cells.Item(56, 2).Value = "KERNEX"
$wsPerf.Cells.Item(56, 3).Value = 21.2368
$wsPerf.Cells.Item(58, 2).Value = "PRIVISCL"
$wsPerf.Cells.Item(58, 3).Value = 20.495
$wsPerf.Cells.Item(59, 2).Value = "INDRAMEDCO"
$wsPerf.Cells.Item(59, 3).Value = 20.4052
$wsPerf.Cells.Item(61, 2).Value = "HINDCOPPER"
$wsPerf.Cells.Item(61, 3).Value = 20.0898
$wsPerf.Cells.Item(62, 2).Value = "BHARATWIRE"
$wsPerf.Cells.Item(62, 3).Value = 20.0711
$wsPerf.Cells.Item(63, 3).Value = 19.9463
$wsPerf.Cells.Item(66, 2).Value = "FEDERALBNK"
$wsPerf.Cells.Item(66, 3).Value = 19.6477
$wsPerf.Cells.Item(67, 2).Value = "BANKINDIA"
$wsPerf.Cells.Item(67, 3).Value = 19.4504
$wsPerf.Cells.Item(68, 3).Value = 19.3027
$wsPerf.Cells.Item(69, 3).Value = 18.9452
$wsPerf.Cells.Item(70, 3).Value = 18.9196
$wsPerf.Cells.Item(73, 3).Value = 18.5523
$wsPerf.Cells.Item(74, 2).Value = "THOMASCOTT"
$wsPerf.Cells.Item(74, 3).Value = 18.4092
$wsPerf.Cells.Item(75, 2).Value = "CEATLTD"
$wsPerf.Cells.Item(75, 3).Value = 18.3784
$wsPerf.Cells.Item(76, 3).Value = 18.3079

# --- distance from Dma50 sheet ---
$wsDma = $wb.Worksheets.Item("distance from Dma50")
$wsDma.Cells.Item(2, 3).Value = 10.157
$wsDma.Cells.Item(3, 3).Value = 7.378
$wsDma.Cells.Item(4, 3).Value = 6.0743
$wsDma.Cells.Item(5, 3).Value = 5.2133
$wsDma.Cells.Item(6, 3).Value = 5.1541
$wsDma.Cells.Item(7, 3).Value = 4.9427
$wsDma.Cells.Item(8, 3).Value = 4.5122
$wsDma.Cells.Item(9, 3).Value = 4.3891
$wsDma.Cells.Item(10, 3).Value = 3.8345
$wsDma.Cells.Item(11, 3).Value = 3.4368
$wsDma.Cells.Item(12, 3).Value = 3.3442
$wsDma.Cells.Item(13, 3).Value = 3.3223
$wsDma.Cells.Item(14, 3).Value = 3.0614
$wsDma.Cells.Item(15, 3).Value = 3.0161
$wsDma.Cells.Item(16, 3).Value = 2.9376
$wsDma.Cells.Item(17, 3).Value = 2.7806
$wsDma.Cells.Item(18, 2).Value = "CNXSMALLCAP"
$wsDma.Cells.Item(18, 3).Value = 2.4887
$wsDma.Cells.Item(19, 2).Value = "NIFTYCPSE"
$wsDma.Cells.Item(19, 3).Value = 2.4664
$wsDma.Cells.Item(20, 2).Value = "NIFTY50VALUE20"
$wsDma.Cells.Item(20, 3).Value = 2.2906
$wsDma.Cells.Item(21, 2).Value = "CNXNIFTYJUNIOR"
$wsDma.Cells.Item(21, 3).Value = 2.2889
$wsDma.Cells.Item(22, 3).Value = 1.5649
$wsDma.Cells.Item(23, 3).Value = 1.3726
$wsDma.Cells.Item(24, 3).Value = 1.3308
$wsDma.Cells.Item(25, 3).Value = 1.1533
$wsDma.Cells.Item(26, 2).Value = "CNXPHARMA"
$wsDma.Cells.Item(26, 3).Value = 1.0709
$wsDma.Cells.Item(27, 2).Value = "NIFTYGROWSECT15"
$wsDma.Cells.Item(27, 3).Value = 1.0608
$wsDma.Cells.Item(28, 3).Value = 0.6418
$wsDma.Cells.Item(29, 3).Value = 0.258
$wsDma.Cells.Item(30, 3).Value = -2.1175

